$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# This workbook has three sheets (Overview, zh-cn, de-de) each tracking the
# handoff status of the same source file. A "Generate Report for handoff" run
# replaced the old source file (fe439221-...md, handed off with hash
# ca67c59f...) with a new one (17a7fd45-...md, handed off with hash
# 4b8cd1234...). The old run had produced a second, failed-handoff row
# (869b81ce-...md / "Handoff transform failed") in addition to the
# .localization-config row; the new run no longer produces that failed row,
# so the .localization-config row moves up to take its place and the old
# last row is dropped.
# ---------------------------------------------------------------------------

function Update-Hyperlink($ws, $addr, $display) {
    foreach ($h in $ws.Hyperlinks) {
        if ($h.Range.Address() -eq $addr) {
            $h.TextToDisplay = $display
            return
        }
    }
}

function Remove-Hyperlink($ws, $addr) {
    foreach ($h in $ws.Hyperlinks) {
        if ($h.Range.Address() -eq $addr) {
            $h.Delete()
            return
        }
    }
}

# ================= Sheet "Overview" =================
$ws = $wb.Worksheets.Item("Overview")

Update-Hyperlink $ws '$A$2' "17a7fd45-f2b8-4a1e-bf91-9439f2fb2369.md"
Update-Hyperlink $ws '$A$3' ".localization-config"
Remove-Hyperlink $ws '$A$4'

$ws.Range("A2").Value = "17a7fd45-f2b8-4a1e-bf91-9439f2fb2369.md"
$ws.Range("A3").Value = ".localization-config"
$ws.Range("B3").Value = "Not to be localized"
$ws.Range("C3").Value = "Not to be localized"

$ws.Rows.Item(4).Delete()

# ================= Sheet "zh-cn" =================
$ws = $wb.Worksheets.Item("zh-cn")

Update-Hyperlink $ws '$A$2' "17a7fd45-f2b8-4a1e-bf91-9439f2fb2369.md"
Update-Hyperlink $ws '$C$2' "17a7fd45-f2b8-4a1e-bf91-9439f2fb2369.4b8cd12340dda38afcc25b814ffb2baafb554da1.zh-cn.xlf"
Update-Hyperlink $ws '$A$3' ".localization-config"
Remove-Hyperlink $ws '$A$4'

$ws.Range("A2").Value = "17a7fd45-f2b8-4a1e-bf91-9439f2fb2369.md"
$ws.Range("C2").Value = "17a7fd45-f2b8-4a1e-bf91-9439f2fb2369.4b8cd12340dda38afcc25b814ffb2baafb554da1.zh-cn.xlf"
$ws.Range("D2").Value = "2016-01-26 12:11:38"
$ws.Range("A3").Value = ".localization-config"
$ws.Range("B3").Value = "Not to be localized"

$ws.Rows.Item(4).Delete()

# ================= Sheet "de-de" =================
$ws = $wb.Worksheets.Item("de-de")

Update-Hyperlink $ws '$A$2' "17a7fd45-f2b8-4a1e-bf91-9439f2fb2369.md"
Update-Hyperlink $ws '$C$2' "17a7fd45-f2b8-4a1e-bf91-9439f2fb2369.4b8cd12340dda38afcc25b814ffb2baafb554da1.de-de.xlf"
Update-Hyperlink $ws '$A$3' ".localization-config"
Remove-Hyperlink $ws '$A$4'

$ws.Range("A2").Value = "17a7fd45-f2b8-4a1e-bf91-9439f2fb2369.md"
$ws.Range("C2").Value = "17a7fd45-f2b8-4a1e-bf91-9439f2fb2369.4b8cd12340dda38afcc25b814ffb2baafb554da1.de-de.xlf"
$ws.Range("D2").Value = "2016-01-26 12:11:48"
$ws.Range("A3").Value = ".localization-config"
$ws.Range("B3").Value = "Not to be localized"

$ws.Rows.Item(4).Delete()
